$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$CellRef, [string]$NewValue)
    $range = $Worksheet.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $NewValue
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "64.755.51"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "3.425.40"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws "D5" "573.75"
$ws.Range("E5").Value = "  -1.29%  "
Set-TextValue $ws "D6" "159.14"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -2.82%  "
$ws.Range("D9").Value = "3.425.46"
$ws.Range("E9").Value = "  -2.22%  "
Set-TextValue $ws "D10" "7.19"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  -2.51%  "
Set-TextValue $ws "D12" "0.439"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").Value = "4.013.00"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("E15").Value = "  -3.89%  "
Set-TextValue $ws "D16" "27.71"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "64.774.61"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "3.446.49"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("E19").Value = "  -2.00%  "
Set-TextValue $ws "D20" "13.86"
$ws.Range("E20").Value = "  -3.14%  "
Set-TextValue $ws "D21" "380.24"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("E22").Value = "  -3.64%  "
Set-TextValue $ws "D23" "0.548"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("E24").Value = "  +0.09%  "
Set-TextValue $ws "D25" "71.94"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("E26").Value = "  -4.79%  "
Set-TextValue $ws "D27" "9.94"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("E28").Value = "  -1.15%  "
Set-TextValue $ws "D29" "1.00"
$ws.Range("E29").Value = "  -0.23%  "
Set-TextValue $ws "D30" "1.46"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  -3.88%  "
Set-TextValue $ws "D32" "2.01"
$ws.Range("E32").Value = "  -2.77%  "
Set-TextValue $ws "D33" "23.23"
$ws.Range("E33").Value = "  -2.25%  "
Set-TextValue $ws "D34" "7.00"
$ws.Range("E34").Value = "  -2.48%  "
Set-TextValue $ws "D35" "1.58"
$ws.Range("E35").Value = "  +1.36%  "
Set-TextValue $ws "D36" "160.87"
$ws.Range("E36").Value = "  -1.08%  "
Set-TextValue $ws "D37" "1.90"
$ws.Range("E37").Value = "  -3.07%  "
$ws.Range("D38").Value = "2.891.18"
Set-TextValue $ws "D39" "0.0751"
$ws.Range("E39").Value = "  -2.93%  "
Set-TextValue $ws "D40" "6.76"
$ws.Range("E40").Value = "  +3.86%  "
Set-TextValue $ws "D41" "26.36"
$ws.Range("E41").Value = "  -3.94%  "
Set-TextValue $ws "D42" "4.55"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  -0.22%  "
Set-TextValue $ws "D44" "0.0316"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("E45").Value = "  -1.60%  "
Set-TextValue $ws "D46" "25.86"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D47" "2.28"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D48" "318.05"
$ws.Range("E48").Value = "  +1.11%  "
Set-TextValue $ws "D49" "1.08"
$ws.Range("E49").Value = "  -3.99%  "
Set-TextValue $ws "D50" "6.51"
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("E51").Value = "  -3.02%  "
